# Radar_metadata.xlsx edit script
# Applies the content changes described by the commit "Updated R code and metadata file":
#  - New data added to row 4 (radar track-table metadata for Kleive/2016)
#  - New shared string "m201608; m201609; m201610"
#  - Selection moved from E12 to F5
#  - Rows 13-15 lose their custom grey "empty-row" formatting (now plain / date style)
#  - Row 15's label cell becomes bold instead of grey
#  - Rows 16-18 labels lose the grey formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: fill in the previously-empty metadata cells ---
$ws.Range("E4").Value = "robinv216_kleive"
$ws.Range("F4").Value = "m201608; m201609; m201610"
$ws.Range("G4").Value = "track"
$ws.Range("I4").Value = (Get-Date -Year 2016 -Month 10 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J4").Value = "trajectory"
$ws.Range("K4").Value = "distance_travelled"
$ws.Range("L4").Value = "airspeed"
$ws.Range("N4").Value = "score"

# --- Rows 13 & 14: drop the grey "customFormat" row style, keep plain date formatting on H/I ---
$ws.Rows.Item(13).ClearFormats()
$ws.Rows.Item(14).ClearFormats()
$ws.Range("H13").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("I13").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("H14").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("I14").NumberFormat = "yyyy\-mm\-dd;@"

# --- Row 15: drop the grey row style; label becomes bold; H/I keep plain date formatting ---
$ws.Rows.Item(15).ClearFormats()
$ws.Range("H15").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("I15").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("E15").Font.Bold = $true

# --- Rows 16-18: drop the grey label formatting (back to default style) ---
$ws.Range("E16").ClearFormats()
$ws.Range("E17").ClearFormats()
$ws.Range("E18").ClearFormats()

# --- Selection moves from E12 to F5 ---
$ws.Range("F5").Select() | Out-Null
